# Improving field descriptions on excel template
#
# This script reproduces, via Excel COM automation, the content/formatting
# changes made to "Excel Template.xlsx":
#   - Column A/B headers get a trailing "*" (required fields)
#   - The "Policy" column description is expanded
#   - A new description is added above the "Teams" column
#   - The "Custom fields" description (above Custom 1) becomes italic
#   - The cursor/selection is left on C10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (column headers): mark required fields with an asterisk ---
$ws.Range("A2").Value = "Application Name*"
$ws.Range("B2").Value = "Business Criticality*"

# --- Row 1 (descriptions above each header) ---

# Clarify that the policy defaults to the one tied to the Business Criticality
$ws.Range("C1").Value = "Case-sensitive policy name (defaults to the one assigned to the Business Criticality)"

# Add a brand-new description above the "Teams" column (M)
$ws.Range("M1").Value = "Comma-delimited list of teams"
$ws.Range("M1").WrapText = $true

# Make the "Custom fields" note (above Custom 1, column P) italic
$ws.Range("P1").Font.Italic = $true

# --- Leave the active selection where the author left it ---
$ws.Range("C10").Select() | Out-Null
